$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("September")

# Update attendance values (0 -> 3) for specific days; the dependent
# totals (columns J and K) are formulas and will recalc automatically.
$ws.Range("F22").Value = 3
$ws.Range("G22").Value = 3

$ws.Range("F28").Value = 3
$ws.Range("G28").Value = 3

$ws.Range("G30").Value = 3
$ws.Range("G31").Value = 3
$ws.Range("G32").Value = 3
$ws.Range("G33").Value = 3

# Update the active cell selection on the sheet.
$ws.Activate()
$ws.Range("G34").Select()

$wb.Save()
